# Applies the "payer_type_primary" table edits:
#   - drop the "Bluecross" row entirely
#   - rename the table header "payer_type_primary" -> "payer_type_primary2"
#   - rename the last data row "Private or Other" -> "Private" and update
#     its three summary-statistic values
# The sibling "payer_type_secondary" table (same row labels/layout, further
# down in the same big table) must be left untouched.
#
# Note: Cell(...).Range.Text includes a trailing cell-mark (CR + BEL), so
# comparisons below use StartsWith rather than exact equality.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nbsp4 = "$([char]160)$([char]160)$([char]160)$([char]160)"

# Locate the "payer_type_primary" header cell (column 1) without assuming a
# fixed row number, so the script is resilient to unrelated changes above it.
$headerRow = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text.StartsWith("payer_type_primary")) {
        $headerRow = $i
        break
    }
}

if ($headerRow -eq -1) {
    Write-Host "ERROR: could not find payer_type_primary header row"
} else {
    # Row directly below the header is "    Bluecross" - delete it outright.
    $bluecrossRow = $headerRow + 1
    $bluecrossText = $t.Cell($bluecrossRow, 1).Range.Text
    if (-not $bluecrossText.StartsWith($nbsp4 + "Bluecross")) {
        Write-Host "WARNING: unexpected row below header:" $bluecrossText
    }
    $t.Rows.Item($bluecrossRow).Delete()

    # Rename the header cell.
    $t.Cell($headerRow, 1).Range.Text = "payer_type_primary2"

    # After the deletion, rows shifted up by one, so the row that used to be
    # "Private or Other" (2 below Bluecross originally) is now 5 rows below
    # the header: Government, Medicaid, Medicare, No Information,
    # No Insurance, Private or Other.
    $privateRow = $headerRow + 6
    $privateLabel = $t.Cell($privateRow, 1).Range.Text
    if (-not $privateLabel.StartsWith($nbsp4 + "Private or Other")) {
        Write-Host "WARNING: unexpected row for Private or Other:" $privateLabel
    }

    $t.Cell($privateRow, 1).Range.Text = $nbsp4 + "Private"
    $t.Cell($privateRow, 3).Range.Text = "99,030 (46%)"
    $t.Cell($privateRow, 4).Range.Text = "86,255 (45%)"
    $t.Cell($privateRow, 5).Range.Text = "12,775 (53%)"

    Write-Host "payer_type_primary table updated."
}
